$wb = $excel.ActiveWorkbook

# --- Belgium: sheet becomes "select all" (A1:XFD1048576), no explicit active cell change needed ---
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
$belgium.Cells.Select()

# --- Duplicate "Spain" to create "Romania", right after Spain ---
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3533/T3549"
$romania.Activate()
$romania.Range("B2:B4").Select()

# --- Duplicate "Romania" to create "Slovakia", right after Romania ---
$romania.Copy($null, $romania)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("B4").Value = "NGC-4306/T3556/T3566"
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Activate()
$slovakia.Range("B2:B4").Select()
